# Commit: Sync attendance_reports, modules_schedules, and assets from main repo
# Update the "Duration" column (G) for every data row on the
# "General_&_Special_Internal_1" sheet to the new value of 900.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 276 }

$ws.Range("G2:G" + $lastRow).Value = 900
